# Append a new data row (row 80) to each of the 4 worksheets, mirroring the
# existing "time / length / id / actual-length / checksum (+ _DEC variants)"
# records. Column G values are huge integers (~1e23) that must be typed in as
# plain digit strings (not "1.2e+23" literal text) so Excel's value-parser
# stores the exact same double without opportunistically switching the cell
# to a scientific NumberFormat (the source rows never carry that style).
$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{
        Name = "FE_LFT_#1"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x20"
        E = "0xf"
        F = 380
        G = "759863127514710945038336"
        H = 288
        I = 15
    },
    @{
        Name = "FE_LFT_#2"
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x2C"
        E = "0xe"
        F = 400
        G = "568432987514711010443264"
        H = 300
        I = 14
    },
    @{
        Name = "FE_PLT_#1"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x60"
        E = "0x3"
        F = 110
        G = "568631262647113970876416"
        H = 96
        I = 3
    },
    @{
        Name = "FE_PLT_#2"
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x5F"
        E = "0x3"
        F = 110
        G = "985046333984776009023488"
        H = 95
        I = 3
    }
)

foreach ($sd in $sheetsData) {
    $ws = $wb.Worksheets.Item($sd.Name)
    $row = 80

    # Column A: timestamp, formatted like the rest of the column.
    $ws.Cells.Item($row, 1).Value = 45866.49649305556
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Columns B-E: comma-joined hex byte strings, kept as text.
    $ws.Cells.Item($row, 2).Value = $sd.B
    $ws.Cells.Item($row, 3).Value = $sd.C
    $ws.Cells.Item($row, 4).Value = $sd.D
    $ws.Cells.Item($row, 5).Value = $sd.E

    # Columns F-I: decoded numeric counterparts.
    $ws.Cells.Item($row, 6).Value = $sd.F
    $ws.Cells.Item($row, 7).Value = $sd.G
    $ws.Cells.Item($row, 8).Value = $sd.H
    $ws.Cells.Item($row, 9).Value = $sd.I
}
